# Auto-generated Excel COM-interop script
# Applies updated market price / profit figures across the ALC, ARM, BSM, CRP, CUL, GSM, LTW and WVR sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1690.4
$ws.Range("I19").Value = 1605.2222
$ws.Range("K19").Value = 1605.2222
$ws.Range("M19").Value = -1430.2222
$ws.Range("H33").Value = 261
$ws.Range("I33").Value = 129.36363
$ws.Range("K33").Value = 129.36363
$ws.Range("M33").Value = 99.63637
$ws.Range("H62").Value = 3954
$ws.Range("I62").Value = 3976.2222
$ws.Range("J62").Value = 3914
$ws.Range("K62").Value = 3976.2222
$ws.Range("L62").Value = 3914
$ws.Range("M62").Value = -3352.2222
$ws.Range("N62").Value = -5162
$ws.Range("H65").Value = 3954
$ws.Range("I65").Value = 3976.2222
$ws.Range("J65").Value = 3914
$ws.Range("K65").Value = 19881.111
$ws.Range("L65").Value = 19570
$ws.Range("M65").Value = -16761.111
$ws.Range("N65").Value = -25810
$ws.Range("H96").Value = 1000
$ws.Range("J96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("N96").ClearContents()
$ws.Range("H113").Value = 21899.334
$ws.Range("I113").Value = 21899
$ws.Range("K113").Value = 21899
$ws.Range("M113").Value = -18645
$ws.Range("H115").Value = 384
$ws.Range("I115").Value = 342.5
$ws.Range("K115").Value = 1027.5
$ws.Range("M115").Value = 539.5
$ws.Range("H135").Value = 726.7
$ws.Range("I135").Value = 728.875
$ws.Range("K135").Value = 6559.875
$ws.Range("M135").Value = -4024.875
$ws.Range("H137").Value = 3465.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3447.1538
$ws.Range("I32").Value = 3109.4167
$ws.Range("K32").Value = 3109.4167
$ws.Range("M32").Value = -2822.4167
$ws.Range("H45").Value = 2453.7144
$ws.Range("I45").Value = 2436.2
$ws.Range("J45").Value = 2497.5
$ws.Range("K45").Value = 2436.2
$ws.Range("L45").Value = 2497.5
$ws.Range("M45").Value = -2059.2
$ws.Range("N45").Value = -3251.5
$ws.Range("H97").Value = 1024.6666
$ws.Range("I97").Value = 909.6
$ws.Range("K97").Value = 909.6
$ws.Range("M97").Value = -413.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 5038.6113
$ws.Range("I86").Value = 4445.636
$ws.Range("J86").Value = 5970.4287
$ws.Range("K86").Value = 4445.636
$ws.Range("L86").Value = 5970.4287
$ws.Range("M86").Value = -3322.636
$ws.Range("N86").Value = -8216.4287
$ws.Range("H89").Value = 5038.6113
$ws.Range("I89").Value = 4445.636
$ws.Range("J89").Value = 5970.4287
$ws.Range("K89").Value = 22228.18
$ws.Range("L89").Value = 29852.1435
$ws.Range("M89").Value = -16612.18
$ws.Range("N89").Value = -41084.14350000001
$ws.Range("H134").Value = 2330.2856
$ws.Range("I134").Value = 2304.8484
$ws.Range("J134").Value = 2750
$ws.Range("K134").Value = 6914.5452
$ws.Range("L134").Value = 8250
$ws.Range("M134").Value = -4379.5452
$ws.Range("N134").Value = -13320

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 861.8570999999999
$ws.Range("I16").Value = 855.0769
$ws.Range("K16").Value = 855.0769
$ws.Range("M16").Value = -568.0769
$ws.Range("H109").Value = 43425
$ws.Range("J109").Value = 43425
$ws.Range("L109").Value = 43425
$ws.Range("N109").Value = -45505
$ws.Range("H113").Value = 861.8570999999999
$ws.Range("I113").Value = 855.0769
$ws.Range("K113").Value = 855.0769
$ws.Range("M113").Value = 1314.9231
$ws.Range("H134").Value = 1757.1724
$ws.Range("I134").Value = 1427.7727
$ws.Range("J134").Value = 2792.4285
$ws.Range("K134").Value = 4283.3181
$ws.Range("L134").Value = 8377.2855
$ws.Range("M134").Value = -1748.3181
$ws.Range("N134").Value = -13447.2855

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 203.66667
$ws.Range("I33").Value = 105.5
$ws.Range("K33").Value = 633
$ws.Range("M33").Value = -350
$ws.Range("H86").Value = 1075
$ws.Range("I86").Value = 987.5
$ws.Range("J86").Value = 1250
$ws.Range("K86").Value = 2962.5
$ws.Range("L86").Value = 3750
$ws.Range("M86").Value = -1776.5
$ws.Range("N86").Value = -6122
$ws.Range("H89").Value = 1075
$ws.Range("I89").Value = 987.5
$ws.Range("J89").Value = 1250
$ws.Range("K89").Value = 8887.5
$ws.Range("L89").Value = 11250
$ws.Range("M89").Value = -2959.5
$ws.Range("N89").Value = -23106

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("L96").ClearContents()
$ws.Range("N96").Value = 0
$ws.Range("H102").Value = 2207.6667
$ws.Range("I102").Value = 2207.6667
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 2207.6667
$ws.Range("L102").Value = 0
$ws.Range("M102").ClearContents()
$ws.Range("N102").Value = -585.6667000000002
$ws.Range("H113").Value = 1982
$ws.Range("I113").Value = 1982
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1982
$ws.Range("L113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = 188
$ws.Range("H126").Value = 15166.444
$ws.Range("I126").Value = 13312.25
$ws.Range("K126").Value = 39936.75
$ws.Range("M126").Value = -37466.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 51041.145
$ws.Range("I7").Value = 51041.145
$ws.Range("K7").Value = 51041.145
$ws.Range("M7").Value = -50929.145
$ws.Range("H93").Value = 1000
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 1000
$ws.Range("K93").Value = 0
$ws.Range("L93").ClearContents()
$ws.Range("M93").Value = 1000
$ws.Range("N93").Value = -3496
$ws.Range("H126").Value = 51041.145
$ws.Range("I126").Value = 51041.145
$ws.Range("K126").Value = 153123.435
$ws.Range("M126").Value = -150653.435

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 477
$ws.Range("I4").Value = 335.33334
$ws.Range("J4").Value = 618.6667
$ws.Range("K4").Value = 335.33334
$ws.Range("L4").Value = 618.6667
$ws.Range("M4").Value = -222.33334
$ws.Range("N4").Value = -844.6667
$ws.Range("H8").Value = 304
$ws.Range("J8").Value = 304
$ws.Range("L8").Value = 304
$ws.Range("N8").Value = -584
$ws.Range("H122").Value = 3742.0908
$ws.Range("I122").Value = 3687.111
$ws.Range("J122").Value = 3989.5
$ws.Range("K122").Value = 11061.333
$ws.Range("L122").Value = 11968.5
$ws.Range("M122").Value = -8611.332999999999
$ws.Range("N122").Value = -16868.5
$ws.Range("H126").Value = 1453.8
$ws.Range("I126").Value = 636
$ws.Range("K126").Value = 1908
$ws.Range("M126").Value = 562
